# Daily automated refresh of EPEX Spot / Gaz / CO2 price sheets.
# Adds "30-aug" column to "Prix Spot" and a new 2025-08-28 row to
# "Gaz" and "CO2".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "Prix Spot": append column BZ (30-aug) mirroring the style
# of the preceding header column (BY) and fill in the 24 hourly
# values underneath it.
# ---------------------------------------------------------------
$spot = $wb.Worksheets.Item("Prix Spot")

$spot.Range("BZ1").Value = "30-aug"
$spot.Range("BY1").Copy()
$spot.Range("BZ1").PasteSpecial(-4122)   # xlPasteFormats

$spotValues = @{
    2  = 67.02
    3  = 53.66
    4  = 58.64
    5  = 81.51
    6  = 37.79
    7  = 25.32
    8  = 63.72
    9  = 54.56
    10 = 70.57
    11 = 32
    12 = 12.82
    13 = 5.23
    14 = 0
    15 = -0.01
    16 = -0.01
    17 = 0
    18 = 2.54
    19 = 7.1
    20 = 18.99
    21 = 57.97
    22 = 69.68
    23 = 49.08
    24 = 45.08
    25 = 40.38
}

foreach ($row in $spotValues.Keys) {
    $spot.Cells.Item($row, 78).Value = $spotValues[$row]   # column 78 = BZ
}

# ---------------------------------------------------------------
# Sheet "Gaz": append row 75 with the new daily quote. The date is
# stored as literal text (matching the existing column-A cells), so
# it is entered with a leading apostrophe to stop Excel from
# re-interpreting it as a date serial.
# ---------------------------------------------------------------
$gaz = $wb.Worksheets.Item("Gaz")
$gaz.Range("A75").Value = "'2025-08-28"
$gaz.Range("B75").Value = 30.55

# ---------------------------------------------------------------
# Sheet "CO2": same new-row pattern.
# ---------------------------------------------------------------
$co2 = $wb.Worksheets.Item("CO2")
$co2.Range("A75").Value = "'2025-08-28"
$co2.Range("B75").Value = 70.95
